# Daily attendance processing - 2025-11-09 03:02:16
#
# Normalises the "Recorded By" column (G): when an entry list contains the
# literal token "System" alongside a real recorder (dnasr281@gmail.com or
# backup@backdoor.com), move "System" so it immediately follows the first
# name in the list instead of trailing at the end.

function Transform-RecordedBy {
    param([string]$val)

    if ($null -eq $val -or $val -eq "") { return $val }

    $rawParts = $val.Split(",")
    $parts = @()
    foreach ($p in $rawParts) {
        $parts += $p.Trim()
    }

    # Only touch lists that actually contain an exact "System" token
    # (case-sensitive - distinct from the lowercase "system" account name).
    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p.Equals("System")) { $hasSystem = $true }
    }
    if (-not $hasSystem) { return $val }

    # Only touch rows recorded by dnasr281@gmail.com or backup@backdoor.com;
    # leave admin@admin.com-only rows untouched.
    $hasTarget = $false
    foreach ($p in $parts) {
        if ($p.Equals("dnasr281@gmail.com") -or $p.Equals("backup@backdoor.com")) { $hasTarget = $true }
    }
    if (-not $hasTarget) { return $val }

    $rest = @()
    foreach ($p in $parts) {
        if (-not $p.Equals("System")) { $rest += $p }
    }

    $newParts = @()
    if ($rest.Count -eq 1) {
        $newParts += "System"
        $newParts += $rest[0]
    } else {
        $newParts += $rest[0]
        $newParts += "System"
        for ($i = 1; $i -lt $rest.Count; $i++) {
            $newParts += $rest[$i]
        }
    }

    return ($newParts -join ", ")
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Value2
    if ($null -eq $current) { continue }

    $updated = Transform-RecordedBy $current
    if (-not $updated.Equals($current)) {
        $cell.Value = $updated
    }
}
